$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 196 - this shifts the existing rows 196:275
# down to 197:276 (values + formatting carried along), matching the
# weekly-refresh pattern where a new observation is prepended and the
# oldest one falls off the bottom (old row 275 -> new row 276).
$ws.Rows.Item(196).Insert()

# Populate the new row 196 with the latest observation.
$ws.Cells.Item(196, 1).Value = 7
$ws.Cells.Item(196, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(196, 3).Value = "Ñuble"
$ws.Cells.Item(196, 4).Value = 44825
$ws.Cells.Item(196, 5).Value = 16
$ws.Cells.Item(196, 6).Value = 100112006
$ws.Cells.Item(196, 7).Value = "Repollo"
$ws.Cells.Item(196, 8).Value = "Crespo record"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 120
$ws.Cells.Item(196, 11).Value = 1400
$ws.Cells.Item(196, 12).Value = 1500
$ws.Cells.Item(196, 13).Value = 1450
$ws.Cells.Item(196, 14).Value = "`$/unidad"
$ws.Cells.Item(196, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(196, 16).Value = 1450
$ws.Cells.Item(196, 17).Value = 1
$ws.Cells.Item(196, 18).Value = "Hortaliza"
